$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    ,@(2, -0.3003816376580256, 0.01536372372024277)
    ,@(3, -0.2022429035621641, -0.1756534860214808)
    ,@(4, -0.4327485595522825, -0.6161972710371395)
    ,@(5, -0.4703665893834659, 0.09538666474221985)
    ,@(6, -0.2422910077715766, -0.110860041128652)
    ,@(7, -0.173266634072466, -0.157357012434773)
    ,@(8, -0.2955884557500205, 0.01647929363112344)
    ,@(9, -0.236576585766175, 0.154616845984686)
    ,@(10, -0.3892055465103567, 0.4135456418784163)
    ,@(11, -0.09000876924643625, 0.05674661383638247)
    ,@(12, -0.2000441421235869, 0.2997334123492504)
    ,@(13, -0.08481858373001307, -0.06733761956078414)
    ,@(14, 0.0691363050739195, -0.2824363911205039)
    ,@(15, 0.05584681469847915, 0.1236540056204973)
    ,@(16, -0.008431221782084357, 0.2982481990021739)
    ,@(17, 0.1162806633408986, -0.2574018958885828)
)

foreach ($item in $values) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
}
